# Move ignition conditions and how to calculate thermal properties to input file; working
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VarPropEta")
$ws.Activate()

# --- Add new "Conserv-4" data row at row 31 (same pattern as existing rows 28-30) ---
$ws.Range("A31").Value() = "Conserv-4"
$ws.Range("B31").Value() = 101
$ws.Range("C31").Value() = 601
$ws.Range("D31").Value() = 400

$ws.Range("E31").HorizontalAlignment = -4152
$ws.Range("E31").Value() = "n/a"

$ws.Range("G31").Value() = 70

$ws.Range("H31").NumberFormat = "0.00E+00"
$ws.Range("H31").Value() = 4890000

$ws.Range("I31").Value() = 20.6

# --- Move the "Conclusions:" / "T>=2500 K for v>0" notes down from D32:D33 to D38:D39 ---
$ws.Range("D32").ClearContents()
$ws.Range("D33").ClearContents()

$ws.Range("D38").Value() = "Conclusions:"
$ws.Range("D39").Value() = "T>=2500 K for v>0"

# --- Scroll the frozen-pane view down and move the active selection to A32 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 12
$win.ScrollColumn = 1
$ws.Range("A32").Select()
